$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "Iniziato parte gestione infringement (fatto/terminare)"
#         -> "parte gestione infringement "
#   a) drop the leading "Iniziato " word
#   b) collapse the trailing " (fatto/terminare)" down to a single space
# ------------------------------------------------------------------
$d.Content.Find.Execute("Iniziato parte gestione ", $true, $false, $false, $false, $false, $true, 1, $false, "parte gestione ", 2) | Out-Null

$d.Content.Find.Execute(" (fatto/terminare)", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# ------------------------------------------------------------------
# Edit 2: paragraph "Gestione dello sblocco dall'infragment lock (lato
#         central) (da terminare)"
#   a) strike-through everything up to (but excluding) the closing
#      parenthesis after "central"
#   b) collapse the trailing " (da terminare)" down to a single space
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Gestione dello sblocco dall*infragment*") {
        $targetPara = $p
        break
    }
}

$paraStart = $targetPara.Range.Start
$paraEnd = $targetPara.Range.End

# locate the start of the run to strike ("Gestione dello sblocco dall...")
$strikeStart = $d.Range($paraStart, $paraEnd)
$strikeStart.Find.Execute("Gestione dello sblocco dall") | Out-Null

# locate the end of the run to strike (the word "central", just before ")")
$strikeEndLocator = $d.Range($paraStart, $paraEnd)
$strikeEndLocator.Find.Execute("central") | Out-Null

$strikeRange = $d.Range($strikeStart.Start, $strikeEndLocator.End)
$strikeRange.Font.StrikeThrough = $true

# collapse " (da terminare)" to a single space, scoped to this paragraph
$tailRange = $d.Range($paraStart, $targetPara.Range.End)
$tailRange.Find.Execute(" (da terminare)", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null
